# "data update may 2nd" - append new daily COVID rows for Egypt (May 1 - May 31,
# 2020) below the existing data, which ran through April 30, 2020 (row 50).
# Full case/death numbers are only available for May 1 and May 2; the rest of
# May only has the date filled in (placeholder rows, as in the source sheet).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Column A uses a date number format (same style as the existing rows, e.g.
# A50). Copy that formatting down to the new date cells first so the new
# cells reuse the existing date style instead of minting a new one.
$ws.Range("A50").Copy() | Out-Null
$ws.Range("A51:A81").PasteSpecial(-4122) | Out-Null  # xlPasteFormats

# Date serials for 2020-05-01 .. 2020-05-31 go in column A for rows 51-81.
$startSerial = 43952
for ($i = 51; $i -le 81; $i++) {
    $ws.Cells.Item($i, 1).Value2 = $startSerial + ($i - 51)
}

# Full data is only available for the first two new days (May 1 and May 2).
$ws.Range("B51").Value2 = 358
$ws.Range("C51").Value2 = 5895
$ws.Range("D51").Value2 = 1460
$ws.Range("E51").Value2 = 406
$ws.Range("F51").Value2 = 14
$ws.Range("G51").Value2 = 1875

$ws.Range("B52").Value2 = 298
$ws.Range("C52").Value2 = 6193
$ws.Range("D52").Value2 = 1522
$ws.Range("E52").Value2 = 415
$ws.Range("F52").Value2 = 9
$ws.Range("G52").Value2 = 1970

# Match the author's final view state: scrolled down and G52 selected.
$ws.Range("G52").Select() | Out-Null
